$d = $word.ActiveDocument

# wdReplaceOne = 1  (only replace the first match - several job entries on this
# resume reuse the same boilerplate phrasing, and the diff only touches the
# very first "WebEquate ... Aug 2022 - Present" entry and summary blurb)
$wdReplaceOne = 1

# 1. Swap "Founder & " / "Full Stack Developer" -> "Full Stack Developer" / " & Founder"
#    (only the first WebEquate entry, Aug 2022 - Present)
$d.Content.Find.Execute("Founder & Full Stack Developer", $true, $false, $false, $false, $false, $true, 1, $false, "Full Stack Developer & Founder", $wdReplaceOne)

# 2. Shorten "Rebuilding websites previously built with PHP, MySQL, and Bootstrap, instead using"
#    paragraph. Stop the match right before "Vercel" so the spellcheck
#    proofErr wrapper around that word is left untouched.
$d.Content.Find.Execute("Rebuilding websites previously built with PHP, MySQL, and Bootstrap, instead using Next.js, React, JavaScript / TypeScript, MongoDB, Tailwind, and ", $true, $false, $false, $false, $false, $true, 1, $false, "Rebuilt websites using Next.js, React, JavaScript / TypeScript, MongoDB, Tailwind, and ", $wdReplaceOne)

# 3. Rewrite summary sentence
$d.Content.Find.Execute("Senior front end software engineer with a background in full stack web development and Agile methodology seeks a fully remote work opportunity.", $true, $false, $false, $false, $false, $true, 1, $false, "Senior front end software engineer and full stack web developer seeks fully remote work opportunities.", $wdReplaceOne)
